$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Canada): update M2_1stDate and M2_LastDate
$ws.Range("E11").Value = 29983
$ws.Range("F11").Value = 45170

# Row 28 (Philippines): update M2_1stDate and M2_LastDate
$ws.Range("E28").Value = 29983
$ws.Range("F28").Value = 45170

# Row 44 (Argentina): update M2_Len and M2_LastDate
$ws.Range("C44").Value = 405
$ws.Range("F44").Value = 45170

# Row 50 (Kazakhstan): update M2_Len and M2_LastDate
$ws.Range("C50").Value = 359
$ws.Range("F50").Value = 45200

# Row 52 (Guatemala): update M2_Len and M2_LastDate
$ws.Range("C52").Value = 346
$ws.Range("F52").Value = 45200
